$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 3867
$ws.Range('K3').Value = 3939
$ws.Range('K4').Value = 802
$ws.Range('K5').Value = 273
$ws.Range('K6').Value = 4474
$ws.Range('K7').Value = 13355

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 113
$ws.Range('K5').Value = 29
$ws.Range('K6').Value = 102
$ws.Range('K7').Value = 387
$ws.Range('K8').Value = 903
$ws.Range('K10').Value = 76
$ws.Range('K11').Value = 267
$ws.Range('K12').Value = 24
$ws.Range('K13').Value = 15
$ws.Range('K18').Value = 92
$ws.Range('K19').Value = 409
$ws.Range('K20').Value = 297
$ws.Range('K21').Value = 37
$ws.Range('K23').Value = 133
$ws.Range('K25').Value = 58
$ws.Range('K29').Value = 706
$ws.Range('K33').Value = 553
$ws.Range('K34').Value = 65
$ws.Range('K36').Value = 165
$ws.Range('K37').Value = 456
$ws.Range('K41').Value = 113
$ws.Range('K42').Value = 473
$ws.Range('K44').Value = 124
$ws.Range('K47').Value = 75
$ws.Range('K48').Value = 172
$ws.Range('K50').Value = 73
$ws.Range('K51').Value = 159
$ws.Range('K52').Value = 366
$ws.Range('K53').Value = 181
$ws.Range('K54').Value = 255
$ws.Range('K55').Value = 150
$ws.Range('K62').Value = 5
$ws.Range('K63').Value = 44
$ws.Range('K65').Value = 311
$ws.Range('K66').Value = 45
$ws.Range('K67').Value = 520
$ws.Range('K73').Value = 124
$ws.Range('K77').Value = 92
$ws.Range('K78').Value = 162
$ws.Range('K79').Value = 347
$ws.Range('K80').Value = 48
$ws.Range('K82').Value = 15
$ws.Range('K83').Value = 281
$ws.Range('K85').Value = 603
$ws.Range('K86').Value = 90
$ws.Range('K88').Value = 151
$ws.Range('K89').Value = 184
$ws.Range('K91').Value = 144
$ws.Range('K94').Value = 164
$ws.Range('K97').Value = 114
$ws.Range('K101').Value = 13355

# Sheet 5: Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K6').Value = 94
$ws.Range('K7').Value = 387

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K6').Value = 99
$ws.Range('K7').Value = 267

# Sheet 7: Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K3').Value = 56
$ws.Range('K7').Value = 184

# Sheet 8: South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K3').Value = 200
$ws.Range('K6').Value = 138
$ws.Range('K7').Value = 603

# Sheet 9: Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K2').Value = 94
$ws.Range('K3').Value = 96
$ws.Range('K7').Value = 366

# Sheet 11: Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K2').Value = 42
$ws.Range('K3').Value = 40
$ws.Range('K7').Value = 181

# Sheet 12: Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 260
$ws.Range('K3').Value = 271
$ws.Range('K4').Value = 50
$ws.Range('K6').Value = 300
$ws.Range('K7').Value = 903

# Sheet 13: South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 102
$ws.Range('K7').Value = 281

# Sheet 14: Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 149
$ws.Range('K3').Value = 209
$ws.Range('K6').Value = 159
$ws.Range('K7').Value = 553

# Sheet 16: Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 122
$ws.Range('K3').Value = 153
$ws.Range('K5').Value = 20
$ws.Range('K6').Value = 138
$ws.Range('K7').Value = 456

# Sheet 17: New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 90
$ws.Range('K6').Value = 126
$ws.Range('K7').Value = 311

# Sheet 21: North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K3').Value = 179
$ws.Range('K7').Value = 520

# Sheet 24: Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K2').Value = 44
$ws.Range('K6').Value = 122
$ws.Range('K7').Value = 255

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 202
$ws.Range('K3').Value = 247
$ws.Range('K6').Value = 200
$ws.Range('K7').Value = 706

# Sheet 26: Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K3').Value = 38
$ws.Range('K4').Value = 24
$ws.Range('K6').Value = 88
$ws.Range('K7').Value = 172

# Sheet 27: Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 133
$ws.Range('K7').Value = 409

# Sheet 28: Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K2').Value = 27
$ws.Range('K7').Value = 124

# Sheet 30: Ashburn
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('K4').Value = 4
$ws.Range('K7').Value = 102

# Sheet 31: Hermosa
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('K2').Value = 38
$ws.Range('K7').Value = 113

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 128
$ws.Range('K3').Value = 154
$ws.Range('K5').Value = 3
$ws.Range('K7').Value = 473

# Sheet 33: Boystown
$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('K3').Value = 5
$ws.Range('K6').Value = 15

# Sheet 34: Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 76

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K3').Value = 34
$ws.Range('K7').Value = 162

# Sheet 36: Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K3').Value = 39
$ws.Range('K7').Value = 150

# Sheet 39: Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K2').Value = 40
$ws.Range('K7').Value = 133

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K3').Value = 67
$ws.Range('K7').Value = 144

# Sheet 41: Chinatown
$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('K3').Value = 13
$ws.Range('K7').Value = 37

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K3').Value = 114
$ws.Range('K7').Value = 347

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K3').Value = 89
$ws.Range('K6').Value = 93
$ws.Range('K7').Value = 297

# Sheet 45: Calumet Heights
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K2').Value = 27
$ws.Range('K7').Value = 92

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 165

# Sheet 50: Garfield Ridge
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K2').Value = 19
$ws.Range('K7').Value = 65

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K3').Value = 30
$ws.Range('K6').Value = 69
$ws.Range('K7').Value = 164

# Sheet 52: East Side
$ws = $wb.Worksheets.Item('East Side')
$ws.Range('K3').Value = 23
$ws.Range('K7').Value = 58

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K6').Value = 24
$ws.Range('K7').Value = 75

# Sheet 56: Lincoln Square
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('K6').Value = 40
$ws.Range('K7').Value = 73

# Sheet 59: North Center
$ws = $wb.Worksheets.Item('North Center')
$ws.Range('K6').Value = 24
$ws.Range('K7').Value = 45

# Sheet 62: Portage Park
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K2').Value = 37
$ws.Range('K3').Value = 30
$ws.Range('K6').Value = 49
$ws.Range('K7').Value = 124

# Sheet 64: Albany Park
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K3').Value = 32
$ws.Range('K7').Value = 113

# Sheet 65: West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K2').Value = 20
$ws.Range('K7').Value = 114

# Sheet 68: United Center
$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K3').Value = 43
$ws.Range('K7').Value = 151

# Sheet 70: Armour Square
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('K3').Value = 8
$ws.Range('K6').Value = 11
$ws.Range('K7').Value = 29

# Sheet 72: Streeterville
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K4').Value = 34
$ws.Range('K7').Value = 90

# Sheet 75: Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K2').Value = 42
$ws.Range('K3').Value = 47
$ws.Range('K6').Value = 53
$ws.Range('K7').Value = 159

# Sheet 83: Sheffield & DePaul
$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range('K2').Value = 3
$ws.Range('K6').Value = 15

# Sheet 84: Riverdale
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K6').Value = 12
$ws.Range('K7').Value = 92

# Sheet 87: Rush & Division
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('K6').Value = 24
$ws.Range('K7').Value = 48

# Sheet 91: Beverly
$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('K2').Value = 9
$ws.Range('K7').Value = 24

# Sheet 98: Museum Campus
$ws = $wb.Worksheets.Item('Museum Campus')
$ws.Range('K6').Value = 4
$ws.Range('K7').Value = 5
